$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct two location rows that were matched to the wrong state/city ---
# Row 150: location should be "Biloxi" (was "Memphis"), Mississippi stays as-is.
$ws.Range("A150").Value = "Biloxi"
# Row 162: location should be "Washington" / "District of Columbia" (was "Columbia" / "North Carolina").
$ws.Range("A162").Value = "Washington"
$ws.Range("B162").Value = "District of Columbia"

# --- Tag each matched row (133-162, skipping the already-annotated 155 & 157)
#     with its data source ("Zillow") and match quality ("good", except the
#     "fixed" case in row 147) in columns J and K ---
$ws.Range("J133").Value = "Zillow"

for ($r = 134; $r -le 162; $r++) {
    if ($r -eq 155 -or $r -eq 157) {
        continue
    }
    if ($r -ne 161) {
        $cell = $ws.Cells.Item($r, 10)
        $cell.Value = "Zillow"
        if ($r -ge 143) {
            # from row 143 on, column J picks up the same (bold/black) style as column A
            $cell.Font.Color = 0
        }
    }
    if ($r -eq 147) {
        $ws.Cells.Item($r, 11).Value = "fixed"
    } else {
        $ws.Cells.Item($r, 11).Value = "good"
    }
}

# --- Scroll the active sheet view to where the new data was added ---
$ws.Application.ActiveWindow.ScrollRow = 129
$ws.Range("M156").Select()
